$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Cells.Item(2, 4) "275.11"
Set-TextValue $ws.Cells.Item(2, 5) "-2.08%"
Set-TextValue $ws.Cells.Item(3, 4) "27.26"
Set-TextValue $ws.Cells.Item(3, 5) "1.22%"
Set-TextValue $ws.Cells.Item(4, 4) "4.770"
Set-TextValue $ws.Cells.Item(4, 5) "-3.53%"
Set-TextValue $ws.Cells.Item(5, 4) "0.06320"
Set-TextValue $ws.Cells.Item(5, 5) "-1.69%"
Set-TextValue $ws.Cells.Item(6, 4) "6.928"
Set-TextValue $ws.Cells.Item(6, 5) "-1.14%"
Set-TextValue $ws.Cells.Item(7, 4) "1.347"
Set-TextValue $ws.Cells.Item(7, 5) "32.50%"
Set-TextValue $ws.Cells.Item(8, 4) "0.8778"
Set-TextValue $ws.Cells.Item(8, 5) "-1.20%"
Set-TextValue $ws.Cells.Item(9, 4) "0.1504"
Set-TextValue $ws.Cells.Item(9, 5) "0.37%"
Set-TextValue $ws.Cells.Item(10, 4) "0.05027"
Set-TextValue $ws.Cells.Item(10, 5) "-4.27%"
Set-TextValue $ws.Cells.Item(11, 4) "0.07513"
Set-TextValue $ws.Cells.Item(11, 5) "1.69%"
Set-TextValue $ws.Cells.Item(12, 4) "0.02916"
Set-TextValue $ws.Cells.Item(12, 5) "-6.68%"
Set-TextValue $ws.Cells.Item(13, 4) "0.08994"
Set-TextValue $ws.Cells.Item(13, 5) "-0.87%"
Set-TextValue $ws.Cells.Item(14, 5) "-0.52%"
Set-TextValue $ws.Cells.Item(15, 4) "0.0006354"
Set-TextValue $ws.Cells.Item(15, 5) "0.64%"
Set-TextValue $ws.Cells.Item(16, 4) "0.005829"
Set-TextValue $ws.Cells.Item(16, 5) "-3.59%"
Set-TextValue $ws.Cells.Item(17, 4) "3.444"
Set-TextValue $ws.Cells.Item(17, 5) "-1.59%"
Set-TextValue $ws.Cells.Item(18, 4) "3.295"
Set-TextValue $ws.Cells.Item(18, 5) "-1.57%"
Set-TextValue $ws.Cells.Item(19, 4) "2.272"
Set-TextValue $ws.Cells.Item(19, 5) "-1.09%"
Set-TextValue $ws.Cells.Item(21, 4) "0.1342"
Set-TextValue $ws.Cells.Item(21, 5) "0.94%"
Set-TextValue $ws.Cells.Item(22, 4) "3.910"
Set-TextValue $ws.Cells.Item(22, 5) "-0.53%"
Set-TextValue $ws.Cells.Item(23, 4) "0.04413"
Set-TextValue $ws.Cells.Item(23, 5) "1.12%"
Set-TextValue $ws.Cells.Item(24, 4) "0.001172"
Set-TextValue $ws.Cells.Item(24, 5) "-0.57%"
Set-TextValue $ws.Cells.Item(25, 4) "0.003849"
Set-TextValue $ws.Cells.Item(25, 5) "4.17%"
Set-TextValue $ws.Cells.Item(26, 4) "0.0001199"
Set-TextValue $ws.Cells.Item(26, 5) "-0.18%"
Set-TextValue $ws.Cells.Item(27, 4) "0.0001932"
Set-TextValue $ws.Cells.Item(27, 5) "13.92%"
Set-TextValue $ws.Cells.Item(40, 4) "0.04115"
Set-TextValue $ws.Cells.Item(40, 5) "-0.15%"
Set-TextValue $ws.Cells.Item(41, 4) "0.006831"
Set-TextValue $ws.Cells.Item(41, 5) "2.35%"
Set-TextValue $ws.Cells.Item(42, 4) "0.1175"
Set-TextValue $ws.Cells.Item(42, 5) "-0.57%"
Set-TextValue $ws.Cells.Item(43, 4) "0.002049"
Set-TextValue $ws.Cells.Item(43, 5) "-13.29%"
Set-TextValue $ws.Cells.Item(44, 4) "0.01157"
Set-TextValue $ws.Cells.Item(44, 5) "-7.97%"
Set-TextValue $ws.Cells.Item(45, 4) "0.00005158"
Set-TextValue $ws.Cells.Item(45, 5) "-2.19%"
Set-TextValue $ws.Cells.Item(47, 4) "0.02295"
Set-TextValue $ws.Cells.Item(47, 5) "2.33%"
